$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 17 de Julio de 2020 a las 07:13"

# Row 15 - Pakistan
$ws.Range("B15").Value = 259999
$ws.Range("C15").Value = 2085
$ws.Range("D15").Value = 183737
$ws.Range("E15").Value = 70787
$ws.Range("G15").Value = 49
$ws.Range("H15").Value = 5475

# Row 72 - Kirguistan
$ws.Range("B72").Value = 13101
$ws.Range("C72").Value = 603
$ws.Range("D72").Value = 3821
$ws.Range("E72").Value = 9107
$ws.Range("G72").Value = 6
$ws.Range("H72").Value = 173

# Row 74 - Australia
$ws.Range("D74").Value = 8114
$ws.Range("E74").Value = 3003

# Row 89 - Haiti
$ws.Range("B89").Value = 6948
$ws.Range("C89").Value = 46
$ws.Range("D89").Value = 3606
$ws.Range("E89").Value = 3197

# Row 104 - Tailandia
$ws.Range("B104").Value = 3239
$ws.Range("C104").Value = 3
$ws.Range("D104").Value = 3096
$ws.Range("E104").Value = 85

# Row 111 - Sri Lanka
$ws.Range("B111").Value = 2687
$ws.Range("E111").Value = 669
